# Update automàtic: dades i banners [2026-02-15 17:20]
# Refreshes DATA_EXTRACCIO timestamps and the latest observation values
# pulled from meteo.cat for each station row in the daily summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# HUMITAT_MITJANA_DIA (column H) stores percentages as literal text (e.g. "62%").
# Force text format first so Excel does not reinterpret them as numeric percentages.
$percentCells = @("H8", "H9", "H11", "H13", "H16", "H17", "H20", "H21", "H24", "H28", "H37", "H39", "H40", "H42")
foreach ($addr in $percentCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("E2").Value = "2026-02-15 17:19:06"

# Row 3
$ws.Range("E3").Value = "2026-02-15 17:19:09"
$ws.Range("K3").Value = "6.9 MJ/m2"
$ws.Range("O3").Value = "-5.6 °C"

# Row 4
$ws.Range("E4").Value = "2026-02-15 17:19:11"
$ws.Range("O4").Value = "7.1 °C"

# Row 5
$ws.Range("E5").Value = "2026-02-15 17:19:14"
$ws.Range("O5").Value = "-5.0 °C"

# Row 6
$ws.Range("E6").Value = "2026-02-15 17:19:16"
$ws.Range("K6").Value = "11.9 MJ/m2"
$ws.Range("O6").Value = "8.3 °C"

# Row 7
$ws.Range("E7").Value = "2026-02-15 17:19:19"
$ws.Range("K7").Value = "12.0 MJ/m2"
$ws.Range("O7").Value = "11.4 °C"

# Row 8
$ws.Range("E8").Value = "2026-02-15 17:19:21"
$ws.Range("H8").Value = "62%"
$ws.Range("J8").Value = "1015.8 hPa"
$ws.Range("K8").Value = "12.1 MJ/m2"

# Row 9
$ws.Range("E9").Value = "2026-02-15 17:19:24"
$ws.Range("H9").Value = "46%"
$ws.Range("K9").Value = "11.8 MJ/m2"
$ws.Range("O9").Value = "11.1 °C"

# Row 10
$ws.Range("E10").Value = "2026-02-15 17:19:26"
$ws.Range("K10").Value = "11.9 MJ/m2"
$ws.Range("O10").Value = "7.4 °C"

# Row 11
$ws.Range("E11").Value = "2026-02-15 17:19:29"
$ws.Range("H11").Value = "38%"

# Row 12
$ws.Range("E12").Value = "2026-02-15 17:19:31"

# Row 13
$ws.Range("E13").Value = "2026-02-15 17:19:33"
$ws.Range("H13").Value = "32%"
$ws.Range("J13").Value = "1015.4 hPa"
$ws.Range("K13").Value = "6.8 MJ/m2"
$ws.Range("O13").Value = "6.3 °C"

# Row 14
$ws.Range("E14").Value = "2026-02-15 17:19:36"

# Row 15
$ws.Range("E15").Value = "2026-02-15 17:19:38"

# Row 16
$ws.Range("E16").Value = "2026-02-15 17:19:41"
$ws.Range("H16").Value = "58%"
$ws.Range("I16").Value = "0.3 mm"
$ws.Range("K16").Value = "9.4 MJ/m2"
$ws.Range("O16").Value = "-2.3 °C"

# Row 17
$ws.Range("E17").Value = "2026-02-15 17:19:43"
$ws.Range("H17").Value = "33%"
$ws.Range("K17").Value = "12.8 MJ/m2"

# Row 18
$ws.Range("E18").Value = "2026-02-15 17:19:46"
$ws.Range("K18").Value = "12.0 MJ/m2"
$ws.Range("O18").Value = "7.0 °C"

# Row 19
$ws.Range("E19").Value = "2026-02-15 17:19:48"
$ws.Range("K19").Value = "11.9 MJ/m2"
$ws.Range("O19").Value = "3.0 °C"

# Row 20
$ws.Range("E20").Value = "2026-02-15 17:19:50"
$ws.Range("H20").Value = "58%"
$ws.Range("K20").Value = "13.0 MJ/m2"
$ws.Range("O20").Value = "-3.1 °C"

# Row 21
$ws.Range("E21").Value = "2026-02-15 17:19:53"
$ws.Range("H21").Value = "35%"
$ws.Range("J21").Value = "1014.9 hPa"
$ws.Range("K21").Value = "10.9 MJ/m2"
$ws.Range("O21").Value = "7.6 °C"

# Row 22
$ws.Range("E22").Value = "2026-02-15 17:19:55"
$ws.Range("K22").Value = "12.6 MJ/m2"
$ws.Range("N22").Value = "-6.5 °C 16:41 TU"

# Row 23
$ws.Range("E23").Value = "2026-02-15 17:19:58"
$ws.Range("K23").Value = "13.0 MJ/m2"
$ws.Range("O23").Value = "-4.0 °C"

# Row 24
$ws.Range("E24").Value = "2026-02-15 17:20:00"
$ws.Range("H24").Value = "68%"
$ws.Range("J24").Value = "1018.2 hPa"
$ws.Range("K24").Value = "11.6 MJ/m2"
$ws.Range("O24").Value = "8.5 °C"

# Row 25
$ws.Range("E25").Value = "2026-02-15 17:20:03"
$ws.Range("O25").Value = "-2.0 °C"

# Row 26
$ws.Range("E26").Value = "2026-02-15 17:20:05"

# Row 27
$ws.Range("E27").Value = "2026-02-15 17:20:07"
$ws.Range("O27").Value = "-0.3 °C"

# Row 28
$ws.Range("E28").Value = "2026-02-15 17:20:10"
$ws.Range("H28").Value = "55%"
$ws.Range("K28").Value = "10.9 MJ/m2"
$ws.Range("O28").Value = "6.3 °C"

# Row 29
$ws.Range("E29").Value = "2026-02-15 17:20:12"
$ws.Range("K29").Value = "12.2 MJ/m2"

# Row 30
$ws.Range("E30").Value = "2026-02-15 17:20:14"

# Row 31
$ws.Range("E31").Value = "2026-02-15 17:20:17"
$ws.Range("J31").Value = "1014.2 hPa"
$ws.Range("O31").Value = "9.6 °C"

# Row 32
$ws.Range("E32").Value = "2026-02-15 17:20:19"
$ws.Range("K32").Value = "9.3 MJ/m2"
$ws.Range("L32").Value = "38.2 km/h - 299º 16:30 TU"
$ws.Range("O32").Value = "3.2 °C"

# Row 33
$ws.Range("E33").Value = "2026-02-15 17:20:22"
$ws.Range("J33").Value = "1015.4 hPa"
$ws.Range("O33").Value = "5.4 °C"

# Row 34
$ws.Range("E34").Value = "2026-02-15 17:20:24"
$ws.Range("K34").Value = "12.2 MJ/m2"
$ws.Range("O34").Value = "0.6 °C"

# Row 35
$ws.Range("E35").Value = "2026-02-15 17:20:27"
$ws.Range("J35").Value = "1019.4 hPa"
$ws.Range("K35").Value = "10.6 MJ/m2"
$ws.Range("O35").Value = "3.8 °C"

# Row 36
$ws.Range("E36").Value = "2026-02-15 17:20:29"
$ws.Range("K36").Value = "10.1 MJ/m2"
$ws.Range("O36").Value = "11.4 °C"

# Row 37
$ws.Range("E37").Value = "2026-02-15 17:20:32"
$ws.Range("H37").Value = "51%"
$ws.Range("O37").Value = "6.0 °C"

# Row 38
$ws.Range("E38").Value = "2026-02-15 17:20:34"
$ws.Range("K38").Value = "12.3 MJ/m2"
$ws.Range("O38").Value = "7.5 °C"

# Row 39
$ws.Range("E39").Value = "2026-02-15 17:20:37"
$ws.Range("H39").Value = "56%"
$ws.Range("K39").Value = "9.9 MJ/m2"
$ws.Range("M39").Value = "0.5 °C 14:42 TU"
$ws.Range("O39").Value = "-3.4 °C"

# Row 40
$ws.Range("E40").Value = "2026-02-15 17:20:39"
$ws.Range("H40").Value = "32%"
$ws.Range("O40").Value = "9.1 °C"

# Row 41
$ws.Range("E41").Value = "2026-02-15 17:20:42"
$ws.Range("O41").Value = "12.1 °C"

# Row 42
$ws.Range("E42").Value = "2026-02-15 17:20:44"
$ws.Range("H42").Value = "52%"
$ws.Range("O42").Value = "10.7 °C"

# Row 43
$ws.Range("E43").Value = "2026-02-15 17:20:47"
$ws.Range("O43").Value = "5.8 °C"

# Row 44
$ws.Range("E44").Value = "2026-02-15 17:20:49"
$ws.Range("K44").Value = "9.7 MJ/m2"

# Row 45
$ws.Range("E45").Value = "2026-02-15 17:20:52"
$ws.Range("J45").Value = "1023.8 hPa"
$ws.Range("K45").Value = "4.5 MJ/m2"
$ws.Range("O45").Value = "0.6 °C"

# Row 46
$ws.Range("E46").Value = "2026-02-15 17:20:54"
$ws.Range("O46").Value = "11.4 °C"
